$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before G (old G -> H)
$ws.Columns("G").Insert()

# Copy number/border formatting from column F into the new column G (rows 5-24)
$ws.Range("F5:F24").Copy()
$ws.Range("G5:G24").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New header / content for the inserted "Thanh toán" column
$ws.Range("G4").Value = "Thanh toán"
$ws.Range("G5").Value = "Đã thanh toán cho Đoàn và Duy ngày 7/9/2019"

# Left-align the new column's data cells, and wrap the text in G5
$ws.Range("G5:G24").HorizontalAlignment = -4131
$ws.Range("G5").WrapText = $true

# Column widths (closest achievable given the engine's column-width quantization)
$ws.Columns("B").ColumnWidth = 32
$ws.Columns("E").ColumnWidth = 12.666666666666666
$ws.Columns("G").ColumnWidth = 27.5
$ws.Columns("H").ColumnWidth = 40.666666666666664

# Update view: clear the fixed top-left cell and move the selection
$ws.Application.ActiveWindow.ScrollRow = 1
[void]$ws.Range("G11").Select()
